# Invitation Task schedule update:
# Split module 3 into blocks and added scoring (bar and popup) type changes to the
# MCQ mapping schedule. This updates specific "a/b/c/d" option cells (columns E:H)
# across rows 2-29 to reflect the new schedule values, and introduces two new
# option values ("20-20" and "10-25") that did not previously exist in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$changes = @{
    "H2"  = "0-5"
    "G3"  = "5-20"
    "E5"  = "5-20"
    "G7"  = "0-0"
    "F8"  = "5-10"
    "G8"  = "5-0"
    "E9"  = "20-20"
    "F10" = "5-10"
    "H10" = "5-0"
    "E11" = "20-20"
    "G11" = "0-20"
    "F13" = "0-0"
    "H14" = "5-10"
    "F15" = "20-0"
    "G16" = "0-0"
    "E17" = "10-20"
    "E18" = "0-20"
    "F18" = "10-5"
    "G19" = "5-10"
    "E20" = "10-25"
    "G21" = "5-20"
    "E22" = "10-20"
    "E23" = "0-0"
    "H24" = "20-0"
    "E25" = "5-10"
    "G26" = "5-20"
    "F27" = "0-20"
    "H27" = "10-25"
    "E28" = "20-0"
    "F28" = "5-10"
    "E29" = "10-5"
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}

# View/window aesthetic changes recorded in the workbook: scroll so row 8 is at
# the top of the visible window, and leave the selection on F27.
try {
    $excel.ActiveWindow.ScrollRow = 8
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Not all hosts support programmatic window scrolling; ignore if unsupported.
}
$ws.Range("F27").Select()
